$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the summary sheet ("总计"): insert the new 2022-Q4 row after the
#    header and shift the existing quarter rows down by one.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Clone the index columns style onto the not-yet-existing A7 before filling it.
$summary.Range("A2").Copy($summary.Range("A7"))

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 30
$summary.Range("D2").Value = 19.69

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 27
$summary.Range("D3").Value = 20.64

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 31
$summary.Range("D4").Value = 22.73

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2021-Q4"
$summary.Range("C5").Value = 26
$summary.Range("D5").Value = 15.67

$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2021-Q3"
$summary.Range("C6").Value = 21
$summary.Range("D6").Value = 10.8

$summary.Range("A7").Value = 5
$summary.Range("B7").Value = "2021-Q2"
$summary.Range("C7").Value = 27
$summary.Range("D7").Value = 16.4

Write-Host "summary sheet (总计) updated"

# ---------------------------------------------------------------------------
# 2. Insert the brand-new "2022-Q4" worksheet right after "总计" (i.e. before
#    the current "2022-Q2" sheet), and populate it with the quarterly fund
#    holdings table.
# ---------------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item(2)
$q4 = $wb.Worksheets.Add($refSheet)
$q4.Name = "2022-Q4"

# Match the header / index-column styling used on the other quarter sheets.
$refSheet.Range("B1:H1").Copy($q4.Range("B1:H1"))
$refSheet.Range("A2").Copy($q4.Range("A2"))

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "'159920"
$q4.Range("C2").Value = "华夏恒生ETF（QDII）"
$q4.Range("D2").Value = "'163.44"
$q4.Range("E2").Value = "'93.95"
$q4.Range("F2").Value = "'7.26"
$q4.Range("G2").Value = "'11.8657"
$q4.Range("H2").Value = 3

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "'513660"
$q4.Range("C3").Value = "华夏沪港通恒生ETF"
$q4.Range("D3").Value = "'30.20"
$q4.Range("E3").Value = "'97.65"
$q4.Range("F3").Value = "'7.63"
$q4.Range("G3").Value = "'2.3043"
$q4.Range("H3").Value = 4

$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "'513550"
$q4.Range("C4").Value = "华泰柏瑞中证港股通50ETF"
$q4.Range("D4").Value = "'25.76"
$q4.Range("E4").Value = "'98.40"
$q4.Range("F4").Value = "'8.59"
$q4.Range("G4").Value = "'2.2128"
$q4.Range("H4").Value = 4

$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "'010365"
$q4.Range("C5").Value = "鹏华港股通中证香港银行投资指数（LOF）C"
$q4.Range("D5").Value = "'6.46"
$q4.Range("E5").Value = "'94.47"
$q4.Range("F5").Value = "'13.80"
$q4.Range("G5").Value = "'0.8915"
$q4.Range("H5").Value = 4

$q4.Range("A6").Value = 4
$q4.Range("B6").Value = "'513600"
$q4.Range("C6").Value = "南方恒生ETF"
$q4.Range("D6").Value = "'7.39"
$q4.Range("E6").Value = "'99.07"
$q4.Range("F6").Value = "'7.90"
$q4.Range("G6").Value = "'0.5838"
$q4.Range("H6").Value = 4

$q4.Range("A7").Value = 5
$q4.Range("B7").Value = "'517300"
$q4.Range("C7").Value = "国寿安保中证沪港深300ETF"
$q4.Range("D7").Value = "'20.94"
$q4.Range("E7").Value = "'99.33"
$q4.Range("F7").Value = "'1.71"
$q4.Range("G7").Value = "'0.3581"
$q4.Range("H7").Value = 8

$q4.Range("A8").Value = 6
$q4.Range("B8").Value = "'501025"
$q4.Range("C8").Value = "鹏华港股通中证香港银行投资指数（LOF）A"
$q4.Range("D8").Value = "'2.38"
$q4.Range("E8").Value = "'94.47"
$q4.Range("F8").Value = "'13.80"
$q4.Range("G8").Value = "'0.3284"
$q4.Range("H8").Value = 4

$q4.Range("A9").Value = 7
$q4.Range("B9").Value = "'513900"
$q4.Range("C9").Value = "华安CES港股通精选100ETF"
$q4.Range("D9").Value = "'1.81"
$q4.Range("E9").Value = "'97.66"
$q4.Range("F9").Value = "'9.17"
$q4.Range("G9").Value = "'0.1660"
$q4.Range("H9").Value = 4

$q4.Range("A10").Value = 8
$q4.Range("B10").Value = "'007354"
$q4.Range("C10").Value = "创金合信港股通量化股票A"
$q4.Range("D10").Value = "'3.27"
$q4.Range("E10").Value = "'90.44"
$q4.Range("F10").Value = "'4.71"
$q4.Range("G10").Value = "'0.1540"
$q4.Range("H10").Value = 4

$q4.Range("A11").Value = 9
$q4.Range("B11").Value = "'517080"
$q4.Range("C11").Value = "汇添富中证沪港深500ETF"
$q4.Range("D11").Value = "'5.13"
$q4.Range("E11").Value = "'92.39"
$q4.Range("F11").Value = "'2.51"
$q4.Range("G11").Value = "'0.1288"
$q4.Range("H11").Value = 5

$q4.Range("A12").Value = 10
$q4.Range("B12").Value = "'517000"
$q4.Range("C12").Value = "银华中证沪港深500ETF"
$q4.Range("D12").Value = "'5.17"
$q4.Range("E12").Value = "'94.61"
$q4.Range("F12").Value = "'2.44"
$q4.Range("G12").Value = "'0.1261"
$q4.Range("H12").Value = 5

$q4.Range("A13").Value = 11
$q4.Range("B13").Value = "'517100"
$q4.Range("C13").Value = "富国中证沪港深500ETF"
$q4.Range("D13").Value = "'3.67"
$q4.Range("E13").Value = "'99.59"
$q4.Range("F13").Value = "'2.57"
$q4.Range("G13").Value = "'0.0943"
$q4.Range("H13").Value = 5

$q4.Range("A14").Value = 12
$q4.Range("B14").Value = "'006810"
$q4.Range("C14").Value = "泰康港股通中证香港银行投资指数C"
$q4.Range("D14").Value = "'0.58"
$q4.Range("E14").Value = "'94.66"
$q4.Range("F14").Value = "'13.99"
$q4.Range("G14").Value = "'0.0811"
$q4.Range("H14").Value = 4

$q4.Range("A15").Value = 13
$q4.Range("B15").Value = "'006809"
$q4.Range("C15").Value = "泰康港股通中证香港银行投资指数A"
$q4.Range("D15").Value = "'0.55"
$q4.Range("E15").Value = "'94.66"
$q4.Range("F15").Value = "'13.99"
$q4.Range("G15").Value = "'0.0769"
$q4.Range("H15").Value = 4

$q4.Range("A16").Value = 14
$q4.Range("B16").Value = "'159712"
$q4.Range("C16").Value = "国泰中证港股通50ETF"
$q4.Range("D16").Value = "'0.74"
$q4.Range("E16").Value = "'97.27"
$q4.Range("F16").Value = "'9.00"
$q4.Range("G16").Value = "'0.0666"
$q4.Range("H16").Value = 1

$q4.Range("A17").Value = 15
$q4.Range("B17").Value = "'513990"
$q4.Range("C17").Value = "招商上证港股通ETF"
$q4.Range("D17").Value = "'0.71"
$q4.Range("E17").Value = "'99.00"
$q4.Range("F17").Value = "'7.29"
$q4.Range("G17").Value = "'0.0518"
$q4.Range("H17").Value = 3

$q4.Range("A18").Value = 16
$q4.Range("B18").Value = "'159711"
$q4.Range("C18").Value = "华夏中证港股通50ETF"
$q4.Range("D18").Value = "'0.35"
$q4.Range("E18").Value = "'96.70"
$q4.Range("F18").Value = "'8.36"
$q4.Range("G18").Value = "'0.0293"
$q4.Range("H18").Value = 4

$q4.Range("A19").Value = 17
$q4.Range("B19").Value = "'012990"
$q4.Range("C19").Value = "天弘国证港股通50指数C"
$q4.Range("D19").Value = "'0.37"
$q4.Range("E19").Value = "'95.03"
$q4.Range("F19").Value = "'7.39"
$q4.Range("G19").Value = "'0.0273"
$q4.Range("H19").Value = 3

$q4.Range("A20").Value = 18
$q4.Range("B20").Value = "'006106"
$q4.Range("C20").Value = "景顺长城量化港股通股票"
$q4.Range("D20").Value = "'0.55"
$q4.Range("E20").Value = "'81.25"
$q4.Range("F20").Value = "'4.39"
$q4.Range("G20").Value = "'0.0241"
$q4.Range("H20").Value = 4

$q4.Range("A21").Value = 19
$q4.Range("B21").Value = "'501309"
$q4.Range("C21").Value = "国泰恒生港股通指数（LOF）"
$q4.Range("D21").Value = "'0.36"
$q4.Range("E21").Value = "'92.87"
$q4.Range("F21").Value = "'6.27"
$q4.Range("G21").Value = "'0.0226"
$q4.Range("H21").Value = 3

$q4.Range("A22").Value = 20
$q4.Range("B22").Value = "'162416"
$q4.Range("C22").Value = "华宝港股通恒生香港35指数（LOF）"
$q4.Range("D22").Value = "'0.25"
$q4.Range("E22").Value = "'90.59"
$q4.Range("F22").Value = "'8.66"
$q4.Range("G22").Value = "'0.0216"
$q4.Range("H22").Value = 3

$q4.Range("A23").Value = 21
$q4.Range("B23").Value = "'005707"
$q4.Range("C23").Value = "富国港股通量化精选股票A"
$q4.Range("D23").Value = "'0.25"
$q4.Range("E23").Value = "'91.95"
$q4.Range("F23").Value = "'7.18"
$q4.Range("G23").Value = "'0.0180"
$q4.Range("H23").Value = 2

$q4.Range("A24").Value = 22
$q4.Range("B24").Value = "'517010"
$q4.Range("C24").Value = "易方达中证沪港深500ETF"
$q4.Range("D24").Value = "'0.48"
$q4.Range("E24").Value = "'93.66"
$q4.Range("F24").Value = "'2.48"
$q4.Range("G24").Value = "'0.0119"
$q4.Range("H24").Value = 5

$q4.Range("A25").Value = 23
$q4.Range("B25").Value = "'007357"
$q4.Range("C25").Value = "创金合信港股通量化股票C"
$q4.Range("D25").Value = "'0.24"
$q4.Range("E25").Value = "'90.44"
$q4.Range("F25").Value = "'4.71"
$q4.Range("G25").Value = "'0.0113"
$q4.Range("H25").Value = 4

$q4.Range("A26").Value = 24
$q4.Range("B26").Value = "'012989"
$q4.Range("C26").Value = "天弘国证港股通50指数A"
$q4.Range("D26").Value = "'0.15"
$q4.Range("E26").Value = "'95.03"
$q4.Range("F26").Value = "'7.39"
$q4.Range("G26").Value = "'0.0111"
$q4.Range("H26").Value = 3

$q4.Range("A27").Value = 25
$q4.Range("B27").Value = "'160925"
$q4.Range("C27").Value = "大成中华沪深港300指数（LOF）A"
$q4.Range("D27").Value = "'0.28"
$q4.Range("E27").Value = "'92.63"
$q4.Range("F27").Value = "'2.77"
$q4.Range("G27").Value = "'0.0078"
$q4.Range("H27").Value = 5

$q4.Range("A28").Value = 26
$q4.Range("B28").Value = "'517170"
$q4.Range("C28").Value = "华夏中证沪港深500ETF"
$q4.Range("D28").Value = "'0.24"
$q4.Range("E28").Value = "'96.90"
$q4.Range("F28").Value = "'2.54"
$q4.Range("G28").Value = "'0.0061"
$q4.Range("H28").Value = 5

$q4.Range("A29").Value = 27
$q4.Range("B29").Value = "'517030"
$q4.Range("C29").Value = "易方达中证沪港深300ETF"
$q4.Range("D29").Value = "'0.38"
$q4.Range("E29").Value = "'93.68"
$q4.Range("F29").Value = "'1.45"
$q4.Range("G29").Value = "'0.0055"
$q4.Range("H29").Value = 8

$q4.Range("A30").Value = 28
$q4.Range("B30").Value = "'008973"
$q4.Range("C30").Value = "大成中华沪深港300指数（LOF）C"
$q4.Range("D30").Value = "'0.03"
$q4.Range("E30").Value = "'92.63"
$q4.Range("F30").Value = "'2.77"
$q4.Range("G30").Value = "'0.0008"
$q4.Range("H30").Value = 5

$q4.Range("A31").Value = 29
$q4.Range("B31").Value = "'014163"
$q4.Range("C31").Value = "富国港股通量化精选股票C"
$q4.Range("D31").Value = "'0.01"
$q4.Range("E31").Value = "'91.95"
$q4.Range("F31").Value = "'7.18"
$q4.Range("G31").Value = "'0.0007"
$q4.Range("H31").Value = 2

Write-Host "2022-Q4 sheet created and populated"

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
  $s = $wb.Worksheets.Item($i)
  Write-Host $i $s.Name
}
